# Swap the "Periodo Mora" (E) and "Valor Mora" (F) values between rows 17
# and 18 on Hoja1. Both rows belong to the same worker (ORLANDO ORDOÑEZ
# GONZALES); the two account-statement periods were in the wrong order and
# are being swapped back.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$e17 = $ws.Range("E17").Value2
$f17 = $ws.Range("F17").Value2
$e18 = $ws.Range("E18").Value2
$f18 = $ws.Range("F18").Value2

$ws.Range("E17").Value2 = $e18
$ws.Range("F17").Value2 = $f18
$ws.Range("E18").Value2 = $e17
$ws.Range("F18").Value2 = $f17
